$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
